# Update countries & provincias Spain
# Applies the COVID-19 dashboard refresh: swaps the displayed country
# names for four row pairs (their underlying stats move with the new
# name), refreshes the numeric stats for a handful of other rows, and
# bumps the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Row, $Country, $Total, $Nuevos, $Activos, $Recuperados, $Criticos, $Muertes) {
    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $Total
    $ws.Cells.Item($Row, 3).Value = $Nuevos
    $ws.Cells.Item($Row, 4).Value = $Activos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 7).Value = $Criticos
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# Row 28 - Catar (name unchanged, stats refreshed)
Set-Row 28 "Catar" 116481 257 113216 3072 0 193

# Row 42 - Rumania
Set-Row 42 "Rumania" 76355 1392 34523 38636 42 3196

# Row 43 - Bielorrusia
Set-Row 43 "Bielorrusia" 70111 161 68256 1223 5 632

# Row 59 - Suiza
Set-Row 59 "Suiza" 39332 306 33900 3434 0 1998

# Row 68/69 - Serbia <-> Nepal swap names, stats follow the new name
Set-Row 68 "Nepal" 30483 838 18214 12132 11 137
Set-Row 69 "Serbia" 30209 0 27908 1612 0 689

# Row 77/78 - Costa de Marfil <-> Bosnia y Herzegovina swap names
Set-Row 77 "Bosnia y Herzegovina" 17396 367 11157 5713 11 526
Set-Row 78 "Costa de Marfil" 17249 0 14611 2526 0 112

# Row 82 - Madagascar
Set-Row 82 "Madagascar" 14218 64 13206 834 1 178

# Row 84/85 - Sudan <-> Senegal swap names
Set-Row 84 "Senegal" 12689 130 8165 4262 1 262
Set-Row 85 "Sudan" 12623 41 6476 5335 0 812

# Row 136/137 - Angola <-> Islandia swap names
Set-Row 136 "Islandia" 2050 10 1920 120 0 10
Set-Row 137 "Angola" 2044 0 742 1209 0 93

# Row 144 - Malta
Set-Row 144 "Malta" 1546 36 854 682 1 10

# Row 158 - Vietnam
Set-Row 158 "Vietnam" 1009 2 545 439 0 25

# Update the "last refreshed" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 21 de Agosto de 2020 a las 13:36"
